$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J11").Value = "Actual 1"
$ws.Range("K11").Value = "Actual 0"

$ws.Range("H12").Value = "Observed 1"
$ws.Range("I12").Value = "Training"

$ws.Range("I13").Value = "Verication"

$ws.Range("H14").Value = "Observed 0"
$ws.Range("I14").Value = "Training"

$ws.Range("I15").Value = "Verication"

$ws.Range("J24").Select()
